$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Despesa"
$ws.Range("B9").Value = "SERVIÇOS"
$ws.Range("C9").Value = 50
$ws.Range("D9").Value = "26/01/2025"
